# Apply updated "想去人数" (column F) counts per the commit diff.
# Each worksheet is addressed by name so the edit is independent of tab order.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 371   # F5: 370 -> 371
$ws.Cells.Item(10, 6).Value = 3037   # F10: 3036 -> 3037
$ws.Cells.Item(11, 6).Value = 1988   # F11: 1984 -> 1988
$ws.Cells.Item(12, 6).Value = 1397   # F12: 1394 -> 1397
$ws.Cells.Item(18, 6).Value = 2151   # F18: 2149 -> 2151
$ws.Cells.Item(19, 6).Value = 1237   # F19: 1236 -> 1237
$ws.Cells.Item(23, 6).Value = 5260   # F23: 5251 -> 5260
$ws.Cells.Item(24, 6).Value = 1018   # F24: 1017 -> 1018
$ws.Cells.Item(25, 6).Value = 91   # F25: 90 -> 91
$ws.Cells.Item(30, 6).Value = 511   # F30: 506 -> 511
$ws.Cells.Item(32, 6).Value = 96   # F32: 95 -> 96
$ws.Cells.Item(33, 6).Value = 2972   # F33: 2971 -> 2972
$ws.Cells.Item(35, 6).Value = 1160   # F35: 1159 -> 1160
$ws.Cells.Item(42, 6).Value = 931   # F42: 930 -> 931
$ws.Cells.Item(47, 6).Value = 1013   # F47: 1012 -> 1013

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 148627   # F5: 148626 -> 148627
$ws.Cells.Item(16, 6).Value = 307   # F16: 305 -> 307
$ws.Cells.Item(26, 6).Value = 524   # F26: 523 -> 524
$ws.Cells.Item(27, 6).Value = 166   # F27: 164 -> 166
$ws.Cells.Item(28, 6).Value = 323   # F28: 322 -> 323
$ws.Cells.Item(37, 6).Value = 123   # F37: 122 -> 123
$ws.Cells.Item(43, 6).Value = 111   # F43: 108 -> 111
$ws.Cells.Item(44, 6).Value = 111   # F44: 108 -> 111

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(7, 6).Value = 914   # F7: 913 -> 914
$ws.Cells.Item(10, 6).Value = 338   # F10: 337 -> 338
$ws.Cells.Item(11, 6).Value = 2575   # F11: 2571 -> 2575
$ws.Cells.Item(12, 6).Value = 152   # F12: 151 -> 152
$ws.Cells.Item(13, 6).Value = 194   # F13: 190 -> 194
$ws.Cells.Item(14, 6).Value = 974   # F14: 973 -> 974

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 914   # F3: 913 -> 914
$ws.Cells.Item(7, 6).Value = 338   # F7: 337 -> 338
$ws.Cells.Item(8, 6).Value = 2575   # F8: 2571 -> 2575
$ws.Cells.Item(9, 6).Value = 371   # F9: 370 -> 371
$ws.Cells.Item(14, 6).Value = 3037   # F14: 3036 -> 3037
$ws.Cells.Item(15, 6).Value = 1988   # F15: 1984 -> 1988
$ws.Cells.Item(16, 6).Value = 1397   # F16: 1394 -> 1397
$ws.Cells.Item(22, 6).Value = 2151   # F22: 2149 -> 2151
$ws.Cells.Item(23, 6).Value = 152   # F23: 151 -> 152
$ws.Cells.Item(24, 6).Value = 1237   # F24: 1236 -> 1237
$ws.Cells.Item(28, 6).Value = 5261   # F28: 5251 -> 5261
$ws.Cells.Item(29, 6).Value = 91   # F29: 90 -> 91
$ws.Cells.Item(31, 6).Value = 323   # F31: 322 -> 323
$ws.Cells.Item(34, 6).Value = 511   # F34: 507 -> 511
$ws.Cells.Item(35, 6).Value = 96   # F35: 95 -> 96
$ws.Cells.Item(36, 6).Value = 2972   # F36: 2971 -> 2972
$ws.Cells.Item(38, 6).Value = 1160   # F38: 1159 -> 1160
$ws.Cells.Item(45, 6).Value = 931   # F45: 930 -> 931
$ws.Cells.Item(49, 6).Value = 1013   # F49: 1012 -> 1013
$ws.Cells.Item(50, 6).Value = 111   # F50: 108 -> 111
